$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows after row 4 to hold rows 5 and 6
$ws.Rows.Item(5).Insert() | Out-Null
$ws.Rows.Item(6).Insert() | Out-Null

# Row 2
$ws.Range("A2").Value = "Pipeline(steps=[('scaler', RobustScaler()),`n                ('selector',`n                 SelectFromModel(estimator=LinearSVC(dual=False, penalty='l1',`n                                                     random_state=42))),`n                ('model',`n                 DecisionTreeClassifier(class_weight='balanced', max_depth=1,`n                                        min_samples_leaf=5, random_state=42))])"
$ws.Range("B2").Value = 0.6571428571428571
$ws.Range("C2").Value = "{'selector': SelectFromModel(estimator=LinearSVC(dual=False, penalty='l1', random_state=42)), 'scaler': RobustScaler(), 'model__min_samples_split': 2, 'model__min_samples_leaf': 5, 'model__max_features': None, 'model__max_depth': 1, 'model__criterion': 'gini', 'model__class_weight': 'balanced'}"
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = "[1 1 0 0 1 0 0 0 0 1 0 1]"
$ws.Range("F2").Value = "[0 0 1 0 0 0 1 0 1 0 0 0]"
$ws.Range("G2").Value = 77
$ws.Range("H2").Value = 0.7775428571428572
$ws.Range("I2").Value = 0.02966702134589264
$ws.Range("J2").Value = 0.5116190476190476
$ws.Range("K2").Value = 0.08681132582019899

# Row 3
$ws.Range("A3").Value = "Pipeline(steps=[('scaler', None), ('selector', None),`n                ('model',`n                 DecisionTreeClassifier(class_weight='balanced',`n                                        criterion='entropy', max_depth=3,`n                                        max_features='log2',`n                                        random_state=42))])"
$ws.Range("B3").Value = 0.5904761904761905
$ws.Range("C3").Value = "{'selector': None, 'scaler': None, 'model__min_samples_split': 2, 'model__min_samples_leaf': 1, 'model__max_features': 'log2', 'model__max_depth': 3, 'model__criterion': 'entropy', 'model__class_weight': 'balanced'}"
$ws.Range("D3").Value = 0.4615384615384615
$ws.Range("E3").Value = "[1 1 0 1 0 0 1 0 1 1 1 0]"
$ws.Range("F3").Value = "[0 1 1 1 1 0 0 1 0 1 0 0]"
$ws.Range("G3").Value = 69
$ws.Range("H3").Value = 0.7762106537530267
$ws.Range("I3").Value = 0.03074252530652464
$ws.Range("J3").Value = 0.4987086359967717
$ws.Range("K3").Value = 0.08070979069907498

# Row 4
$ws.Range("A4").Value = "Pipeline(steps=[('scaler', StandardScaler()), ('selector', None),`n                ('model',`n                 DecisionTreeClassifier(class_weight='balanced',`n                                        criterion='entropy', max_depth=10,`n                                        min_samples_leaf=5, min_samples_split=5,`n                                        random_state=42))])"
$ws.Range("B4").Value = 0.5904761904761904
$ws.Range("C4").Value = "{'selector': None, 'scaler': StandardScaler(), 'model__min_samples_split': 5, 'model__min_samples_leaf': 5, 'model__max_features': None, 'model__max_depth': 10, 'model__criterion': 'entropy', 'model__class_weight': 'balanced'}"
$ws.Range("D4").Value = 0.3636363636363636
$ws.Range("E4").Value = "[1 0 1 1 1 1 0 1 0 1 0 1]"
$ws.Range("F4").Value = "[1 0 0 0 1 0 0 0 0 0 1 0]"
$ws.Range("G4").Value = 42
$ws.Range("H4").Value = 0.777751756440281
$ws.Range("I4").Value = 0.0328583136139833
$ws.Range("J4").Value = 0.5081967213114754
$ws.Range("K4").Value = 0.08314903700646621

# Row 5
$ws.Range("A5").Value = "Pipeline(steps=[('scaler', None),`n                ('selector',`n                 SelectFromModel(estimator=ExtraTreesClassifier(random_state=42))),`n                ('model',`n                 DecisionTreeClassifier(class_weight='balanced',`n                                        criterion='entropy', max_depth=15,`n                                        max_features='sqrt', min_samples_leaf=9,`n                                        min_samples_split=9,`n                                        random_state=42))])"
$ws.Range("B5").Value = 0.6285714285714286
$ws.Range("C5").Value = "{'selector': SelectFromModel(estimator=ExtraTreesClassifier(random_state=42)), 'scaler': None, 'model__min_samples_split': 9, 'model__min_samples_leaf': 9, 'model__max_features': 'sqrt', 'model__max_depth': 15, 'model__criterion': 'entropy', 'model__class_weight': 'balanced'}"
$ws.Range("D5").Value = 0.75
$ws.Range("E5").Value = "[1 1 0 0 0 0 1 0 1 1 1 1]"
$ws.Range("F5").Value = "[0 1 0 1 0 1 1 1 1 1 1 1]"
$ws.Range("G5").Value = 11
$ws.Range("H5").Value = 0.7596221959858322
$ws.Range("I5").Value = 0.03289739517642348
$ws.Range("J5").Value = 0.512396694214876
$ws.Range("K5").Value = 0.07365827583245997

# Row 6
$ws.Range("A6").Value = "Pipeline(steps=[('scaler', None), ('selector', None),`n                ('model',`n                 DecisionTreeClassifier(class_weight='balanced', max_depth=15,`n                                        max_features='log2', min_samples_leaf=9,`n                                        min_samples_split=7,`n                                        random_state=42))])"
$ws.Range("B6").Value = 0.6190476190476191
$ws.Range("C6").Value = "{'selector': None, 'scaler': None, 'model__min_samples_split': 7, 'model__min_samples_leaf': 9, 'model__max_features': 'log2', 'model__max_depth': 15, 'model__criterion': 'gini', 'model__class_weight': 'balanced'}"
$ws.Range("D6").Value = 0.5714285714285715
$ws.Range("E6").Value = "[1 1 1 1 0 0 0 0 1 1 0 0]"
$ws.Range("F6").Value = "[1 1 0 0 0 1 1 1 1 1 1 0]"
$ws.Range("G6").Value = 14
$ws.Range("H6").Value = 0.7644613583138172
$ws.Range("I6").Value = 0.03803203021572739
$ws.Range("J6").Value = 0.5025761124121779
$ws.Range("K6").Value = 0.08522071919017508
